$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(669, 4).Value = 2.818060649304283
$ws.Cells.Item(669, 5).Value = 77.11508748397709
$ws.Cells.Item(671, 4).Value = 6.214969801346295
$ws.Cells.Item(671, 5).Value = 97.67125372082616
$ws.Cells.Item(672, 4).Value = 3.579151351801821
$ws.Cells.Item(672, 5).Value = 77.00568462032463
$ws.Cells.Item(673, 5).Value = 51.51551410311617
$ws.Cells.Item(675, 4).Value = 18.80272513372305
$ws.Cells.Item(675, 5).Value = 95.69973602207925
$ws.Cells.Item(676, 4).Value = 4.940232790045685
$ws.Cells.Item(676, 5).Value = 86.29492790210902
$ws.Cells.Item(677, 4).Value = 23.5340070274588
$ws.Cells.Item(677, 5).Value = 99.5580348881647
$ws.Cells.Item(679, 4).Value = 48.45987320770406
$ws.Cells.Item(679, 5).Value = 97.74792125434936
$ws.Cells.Item(680, 4).Value = 3.258671848892099
$ws.Cells.Item(680, 5).Value = 19.45088198724385
$ws.Cells.Item(682, 4).Value = 51.76696834115641
$ws.Cells.Item(682, 5).Value = 99.33284854808117
$ws.Cells.Item(683, 4).Value = 5.32645116992087
$ws.Cells.Item(683, 5).Value = 38.57067583297124
$ws.Cells.Item(684, 4).Value = 41.30339991193689
$ws.Cells.Item(684, 5).Value = 87.75981356705135
$ws.Cells.Item(686, 4).Value = 3.669690689782173
$ws.Cells.Item(686, 5).Value = 95.58451957564684
$ws.Cells.Item(687, 5).Value = 4.600091855124408
$ws.Cells.Item(688, 5).Value = 5.394627728385085
$ws.Cells.Item(689, 4).Value = 27.67852306455712
$ws.Cells.Item(689, 5).Value = 99.1533203837518
$ws.Cells.Item(691, 4).Value = 0.3618983406209854
$ws.Cells.Item(691, 5).Value = 55.6224244379651
$ws.Cells.Item(692, 4).Value = 14.84082921730442
$ws.Cells.Item(692, 5).Value = 68.97245620203407
$ws.Cells.Item(693, 4).Value = 18.79938316501084
$ws.Cells.Item(693, 5).Value = 43.37428371653264
$ws.Cells.Item(694, 4).Value = 45.91688908787116
$ws.Cells.Item(694, 5).Value = 91.02804851772109
$ws.Cells.Item(695, 4).Value = 2.875622778370469
$ws.Cells.Item(695, 5).Value = 36.95392077248341
$ws.Cells.Item(696, 4).Value = 43.99979564942866
$ws.Cells.Item(696, 5).Value = 96.90333295516733
$ws.Cells.Item(697, 4).Value = 16.59203973931322
$ws.Cells.Item(697, 5).Value = 64.01675974241402
$ws.Cells.Item(814, 4).Value = 21.2269512945991
$ws.Cells.Item(814, 5).Value = 91.9875912620228
$ws.Cells.Item(816, 4).Value = 2.78261046459123
$ws.Cells.Item(816, 5).Value = 52.42770943367551
$ws.Cells.Item(817, 4).Value = 0.3256097310183747
$ws.Cells.Item(817, 5).Value = 22.45613745420824
$ws.Cells.Item(818, 4).Value = 0.4059546626854059
$ws.Cells.Item(818, 5).Value = 71.7078406104323
$ws.Cells.Item(820, 4).Value = 16.00621241445315
$ws.Cells.Item(820, 5).Value = 94.22431697750385
$ws.Cells.Item(821, 4).Value = 6.70592911681136
$ws.Cells.Item(821, 5).Value = 61.55868370840217
$ws.Cells.Item(822, 4).Value = 3.890319539492073
$ws.Cells.Item(822, 5).Value = 71.36515072692058
$ws.Cells.Item(827, 4).Value = 73.07577518868285
$ws.Cells.Item(827, 5).Value = 91.05383239735697
$ws.Cells.Item(828, 4).Value = 6.579839820496411
$ws.Cells.Item(828, 5).Value = 29.85936988989425
$ws.Cells.Item(829, 4).Value = 67.3671885005934
$ws.Cells.Item(829, 5).Value = 99.04618701183504
$ws.Cells.Item(831, 4).Value = 10.07719779551353
$ws.Cells.Item(831, 5).Value = 68.5422885791731
$ws.Cells.Item(832, 5).Value = 0.7602970879352262
$ws.Cells.Item(833, 5).Value = 13.39785699380517
$ws.Cells.Item(834, 4).Value = 46.34282331058927
$ws.Cells.Item(834, 5).Value = 99.3620712443515
$ws.Cells.Item(836, 4).Value = 0.3405114414611988
$ws.Cells.Item(836, 5).Value = 41.55262503110227
$ws.Cells.Item(837, 4).Value = 22.84944815686427
$ws.Cells.Item(837, 5).Value = 94.59667495168202
$ws.Cells.Item(838, 4).Value = 22.48240824770966
$ws.Cells.Item(838, 5).Value = 41.54402398920677
$ws.Cells.Item(841, 4).Value = 39.9575375128921
$ws.Cells.Item(841, 5).Value = 99.43213166678704
$ws.Cells.Item(842, 4).Value = 21.79574534937334
$ws.Cells.Item(842, 5).Value = 76.76700554708007
$ws.Cells.Item(1017, 4).Value = 0.7675438596491228
$ws.Cells.Item(1017, 5).Value = 88.12842825916609
$ws.Cells.Item(1019, 4).Value = 4.451038575667656
$ws.Cells.Item(1020, 4).Value = 0.5649717514124294
$ws.Cells.Item(1021, 4).Value = 0.3891050583657588
$ws.Cells.Item(1023, 4).Value = 1.193576388888889
$ws.Cells.Item(1024, 4).Value = 21.54471544715447
$ws.Cells.Item(1024, 5).Value = 95.17878234725652
$ws.Cells.Item(1025, 4).Value = 9.495548961424333
$ws.Cells.Item(1025, 5).Value = 97.55560145975349
$ws.Cells.Item(1031, 4).Value = 2.390581253439901
$ws.Cells.Item(1031, 5).Value = 45.97910291401318
$ws.Cells.Item(1032, 4).Value = 20
$ws.Cells.Item(1032, 5).Value = 77.36185383244207
$ws.Cells.Item(1034, 4).Value = 3.560830860534125
$ws.Cells.Item(1037, 5).Value = 98.51632047477746
$ws.Cells.Item(1039, 4).Value = 0.525866035627424
$ws.Cells.Item(1039, 5).Value = 73.82960162838035
$ws.Cells.Item(1040, 4).Value = 40.5940594059406
$ws.Cells.Item(1040, 5).Value = 96.67194928684627
$ws.Cells.Item(1041, 4).Value = 3.264095074875772
$ws.Cells.Item(1044, 4).Value = 41.07373868046572
$ws.Cells.Item(1044, 5).Value = 98.15270935960592
$ws.Cells.Item(1045, 4).Value = 4.616347905282332
$ws.Cells.Item(1045, 5).Value = 39.0625
